$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "1 tRNA estimate" block (rows 11-16, columns A-H) ---
$ws.Range("A11").Value = "1 tRNA extimate"
$ws.Range("B11").Value = "Ribosomes"
$ws.Range("B11").NumberFormat = "0.00E+00"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 12
$ws.Range("E11").Formula = "=(4/3)*3.1415*(D11*0.001)^3"
$ws.Range("E11").NumberFormat = "0.00E+00"
$ws.Range("F11").Formula = "=E11*C11"
$ws.Range("F11").NumberFormat = "0.00E+00"
$ws.Range("G11").Formula = "=C11/1621"
$ws.Range("H11").Formula = "=((1.38*10^-23*310)/(6*3.1415*(6.9*10^-4)*(D11*10^-9)))*10^12"
$ws.Range("I11").ClearContents()
$ws.Range("L11").Value = "sub-vol="
$ws.Range("M11").Formula = "=C5/47"
$ws.Range("Q11").Value = "kT/6pietta*a"
$ws.Range("B12").Value = "EF-Tu-tRNA"
$ws.Range("B12").NumberFormat = "0.00E+00"
$ws.Range("C12").Value = 7
$ws.Range("D12").Value = 7.5
$ws.Range("E12").Formula = "=(4/3)*3.1415*(D12*0.001)^3"
$ws.Range("E12").NumberFormat = "0.00E+00"
$ws.Range("F12").Formula = "=E12*C12"
$ws.Range("F12").NumberFormat = "0.00E+00"
$ws.Range("G12").Formula = "=C12/1621"
$ws.Range("H12").Formula = "=((1.38*10^-23*310)/(6*3.1415*(6.9*10^-4)*(D12*10^-9)))*10^12"
$ws.Range("I12").ClearContents()
$ws.Range("L12").Formula = "=C4/M11"
$ws.Range("M12").Value = "dynamic visc water@37C=69*10^-4 kg/(m*s^2)"
$ws.Range("B13").Value = "tRNA synthetase"
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 2.5
$ws.Range("E13").Formula = "=(4/3)*3.1415*(D13*0.001)^3"
$ws.Range("E13").NumberFormat = "0.00E+00"
$ws.Range("F13").Formula = "=E13*C13"
$ws.Range("F13").NumberFormat = "0.00E+00"
$ws.Range("G13").Formula = "=C13/1621"
$ws.Range("H13").Formula = "=((1.38*10^-23*310)/(6*3.1415*(6.9*10^-4)*(D13*10^-9)))*10^12"
$ws.Range("I13").ClearContents()
$ws.Range("L13").Formula = "=C5/M11"
$ws.Range("M13").Value = "k= 1.380*10^-23 kg*m^2/(s^2*kelvin)"
$ws.Range("B14").NumberFormat = "0.00E+00"
$ws.Range("C14").Value = 0
$ws.Range("E14").NumberFormat = "0.00E+00"
$ws.Range("F14").NumberFormat = "0.00E+00"
$ws.Range("G14").Formula = "=C14/1621"
$ws.Range("H14").ClearContents()
$ws.Range("I14").ClearContents()
$ws.Range("L14").Formula = "=C6/M11"
$ws.Range("M14").Value = "Temp = 310 kelvin"
$ws.Range("B15").Value = "aa "
$ws.Range("B15").NumberFormat = "0.00E+00"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0.31
$ws.Range("E15").Formula = "=(4/3)*3.1415*(D15*0.001)^3"
$ws.Range("E15").NumberFormat = "0.00E+00"
$ws.Range("F15").Formula = "=E15*C15"
$ws.Range("F15").NumberFormat = "0.00E+00"
$ws.Range("G15").Formula = "=C15/1621"
$ws.Range("H15").Formula = "=((1.38*10^-23*310)/(6*3.1415*(6.9*10^-4)*(D15*10^-9)))*10^12"
$ws.Range("B16").NumberFormat = "0.00E+00"
$ws.Range("E16").NumberFormat = "0.00E+00"
$ws.Range("F16").Formula = "=SUM(F11:F15)"
$ws.Range("F16").NumberFormat = "0.00E+00"
$ws.Range("H16").ClearContents()
$ws.Range("L16").Formula = "=C8/M11"

# --- Update selection to reflect final cursor position ---
$ws.Range("F13").Select()
